$wb = $excel.ActiveWorkbook
$s1 = $wb.Worksheets.Item(1)
$s2 = $wb.Worksheets.Item(2)
$s3 = $wb.Worksheets.Item(3)

# --- Sheet1: data currently at A1:C2 moves to G6:L7 (non-contiguous columns),
# with a new bold header row added at row 5 (G5/J5/L5).
$s1.Range("G6").Value = $s1.Range("A1").Value2
$s1.Range("J6").Value = $s1.Range("B1").Value2
$s1.Range("L6").Value = $s1.Range("C1").Value2
$s1.Range("G7").Value = $s1.Range("A2").Value2
$s1.Range("J7").Value = $s1.Range("B2").Value2
$s1.Range("L7").Value = $s1.Range("C2").Value2
$s1.Range("A1:C2").Clear()

$s1.Range("G5").Value = "section"
$s1.Range("G5").Font.Bold = $true
$s1.Range("J5").Value = "country"
$s1.Range("J5").Font.Bold = $true
$s1.Range("L5").Value = "award"
$s1.Range("L5").Font.Bold = $true

# --- Sheet2: data currently at A1:C2 shifts down+right to B3:D4 (contiguous
# move), with a new bold header row inserted above it at B2:D2.
$s2.Range("A1:C2").Cut($s2.Range("B3"))

$s2.Range("B2").Value = "section"
$s2.Range("B2").Font.Bold = $true
$s2.Range("C2").Value = "country"
$s2.Range("C2").Font.Bold = $true
$s2.Range("D2").Value = "award"
$s2.Range("D2").Font.Bold = $true

# --- Sheet3: same shift as Sheet2.
$s3.Range("A1:C2").Cut($s3.Range("B3"))

$s3.Range("B2").Value = "section"
$s3.Range("B2").Font.Bold = $true
$s3.Range("C2").Value = "country"
$s3.Range("C2").Font.Bold = $true
$s3.Range("D2").Value = "award"
$s3.Range("D2").Font.Bold = $true

# --- Selections / active sheet: Sheet1 & Sheet2 selections updated but not
# activated; Sheet3 ends up the active/selected tab (matches activeTab="2").
$s1.Range("D10").Select()
$s2.Range("B2").Select()
$s3.Range("B2").Select()
